$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New row 33 : Sno 32, date 2022-07-01 (serial 44743), Research paper entry
# ---------------------------------------------------------------------------
$ws.Range("A33").Value = 32
$ws.Range("B33").Value = 44743
$ws.Range("B33").NumberFormat = "m/d/yy"
$ws.Range("C33").Value = 0.5
$ws.Range("C33").NumberFormat = "h:mm AM/PM"
$ws.Range("D33").Value = 0.5625
$ws.Range("D33").NumberFormat = "h:mm AM/PM"
$ws.Range("E33").Formula = "=D33-C33"
$ws.Range("E33").NumberFormat = "h:mm"
$ws.Range("F33").Value = "Research paper"
$ws.Range("G33").Value = "1. deeplabv3 + paper review and notes"
$ws.Range("G33").WrapText = $true

# ---------------------------------------------------------------------------
# New row 34 : Sno 33, date 2022-07-01, Code entry (deeplabv3+ starter)
# ---------------------------------------------------------------------------
$ws.Range("A34").Value = 33
$ws.Range("B34").Value = 44743
$ws.Range("B34").NumberFormat = "m/d/yy"
$ws.Range("C34").Value = 0.59375
$ws.Range("C34").NumberFormat = "h:mm AM/PM"
$ws.Range("D34").Value = 0.63541666666666663
$ws.Range("D34").NumberFormat = "h:mm AM/PM"
$ws.Range("E34").Formula = "=D34-C34"
$ws.Range("E34").NumberFormat = "h:mm"
$ws.Range("F34").Value = "Code"
$ws.Range("G34").Value = "1. deeplabv3+_starter nb completed"
$ws.Range("G34").WrapText = $true

# ---------------------------------------------------------------------------
# New row 35 : Sno 34, date 2022-07-01, Code entry (deeplabv3+ os8 dice loss)
# ---------------------------------------------------------------------------
$ws.Range("A35").Value = 34
$ws.Range("B35").Value = 44743
$ws.Range("B35").NumberFormat = "m/d/yy"
$ws.Range("C35").Value = 0.70833333333333337
$ws.Range("C35").NumberFormat = "h:mm AM/PM"
$ws.Range("D35").Value = 0.77083333333333337
$ws.Range("D35").NumberFormat = "h:mm AM/PM"
$ws.Range("E35").Formula = "=D35-C35"
$ws.Range("E35").NumberFormat = "h:mm"
$ws.Range("F35").Value = "Code"
$ws.Range("G35").Value = "1. deeplabv3+ os8 version 10ep with dice loss, 1cycle LR 3e-4"
$ws.Range("G35").WrapText = $true

# ---------------------------------------------------------------------------
# Move the "Total Hours" summary row from row 33 down to row 38, leaving
# rows 36-37 blank, and extend the SUM range to cover the new rows.
# ---------------------------------------------------------------------------
$ws.Range("C38").Value = "Total Hours"
$ws.Range("E38").Formula = "=SUM(E2:E37)"
$ws.Range("E38").NumberFormat = "[hh]:mm"

# ---------------------------------------------------------------------------
# Update sheet view: scroll position and active selection
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 28
$win.ScrollColumn = 1
$ws.Range("D36").Select()
